$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update T1 (Checked companies count) ---
$ws.Range("T1").Value = 12

# --- Update existing row 2 ---
$ws.Range("K2").Value = 16.75
$ws.Range("L2").Value = "manage, job, team, process, resource, capacity, match, role, candidates, plan"

# --- Update existing row 3 ---
$ws.Range("A3").Value = "ClearCaptions, LLC"
$ws.Range("K3").Value = 20.82
$ws.Range("L3").Value = "data, com, business, work, res, able, skill, analysis, team, unit"

# --- Helper data for the new rows (4 through 13) ---
$newRows = @(
    @{ A="Ride Health";         B="Workforce Analyst - Fully Remote";                                        K=44.67;  L="com, age, work, health, ride, time, workforce, manage, any, skill" },
    @{ A="UNFI";                 B="FP & A Analyst II- Remote";                                                K=55.68;  L="business, financial, required, experience, remote, work, able, team, unfi, office" },
    @{ A="TieTalent";            B="Business Metrics/Analytics - Remote (Work 8am - 5pm PST Timings)";         K=65.73;  L="data, management, experience, business, work, metrics, skills, ability, portfolio, resource" },
    @{ A="Kforce Inc";           B="Finance Manager, Customer Finance - (Remote)";                             K=69.63;  L="age, finance, service, financial, kforce, customer, pay, team, act, eligible" },
    @{ A="Centene Corporation";  B="Capacity Planning Analyst II";                                              K=85.34;  L="per, act, capacity, planning, work, center, contact, perform, experience, standards" },
    @{ A="Jobgether";            B="Workforce Analyst - (Remote - US)";                                        K=95.55;  L="work, per, staffing, job, workforce, team, time, match, support, candidates" },
    @{ A="Jobgether";            B="Financial Planning & Analysis Senior Analyst - (Remote - USA)";            K=146.78; L="financial, plan, analysis, planning, job, opportunities, performance, match, candidates, business" },
    @{ A="Ryder System, Inc.";   B="Finance Manager - REMOTE";                                                 K=159.8;  L="com, financial, age, ryder, manage, app, work, plan, view, job" },
    @{ A="Ryder System, Inc.";   B="Finance Manager - REMOTE";                                                 K=162.05; L="com, financial, age, ryder, manage, app, work, plan, view, job" },
    @{ A="Ryder System, Inc.";   B="Finance Manager - REMOTE";                                                 K=164.29; L="com, financial, age, ryder, manage, app, work, plan, view, job" }
)

$row = 4
foreach ($item in $newRows) {
    $ws.Cells.Item($row, 1).Value = $item.A
    $ws.Cells.Item($row, 2).Value = $item.B
    $ws.Cells.Item($row, 3).Value = $false
    $ws.Cells.Item($row, 4).Value = $false
    $ws.Cells.Item($row, 5).Value = $false
    $ws.Cells.Item($row, 6).Value = $true
    $ws.Cells.Item($row, 7).Value = $false
    $ws.Cells.Item($row, 8).Value = $true
    $ws.Cells.Item($row, 9).Value = $false
    $ws.Cells.Item($row, 10).Value = "https://www.linkedin.com/jobs/search/"
    $ws.Cells.Item($row, 11).Value = $item.K
    $ws.Cells.Item($row, 12).Value = $item.L
    $row++
}
